$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile x100k ab.)
$data = @(
    @(44313, 2, 9, 144.2076590290018),
    @(44314, 0, 9, 144.2076590290018),
    @(44315, 3, 11, 176.253805479891),
    @(44316, 1, 12, 192.2768787053357),
    @(44317, 1, 11, 176.253805479891),
    @(44318, 1, 9, 144.2076590290018)
)

$startRow = 239
$templateRow = 238

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy formatting from the last existing row (A238) down to the new row's A cell
    $ws.Range("A$templateRow").Copy($ws.Range("A$r"))

    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}
